$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the placeholder link for the Skills Imperative 2035 data source
$ws.Range("B14").Value = "<a href='https://www.gov.uk/government/publications/labour-market-and-skills-projections-2020-to-2035'>Skills Imperative 2035</a>"
$ws.Range("C14").Value = "2035 (16/03/23)"

# Move the active selection to A12
$ws.Range("A12").Select()
